# This script applies the row-data permutation + F-column renumbering
# described by the commit: "new input files generation for english version;
# make only 20 different versions and duplicate many times for 1000 subjects"
#
# Each data row (2..41) keeps its own subject_id/task/block_total/block_scene/trial_block
# (columns A-E) but gets new trial_total (F) and new target_cat/category/cond_cat/
# correct_answer/stimulus/conceptual/perceptual/typicality/n/p_*/r_* values (G,H,I,K,L,M,N,O,P,Q,R,S,T,U,V)
# taken from the permuted source rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,82,"kitchens","kitchens","target","j","stimuli/img_7wul8.png",43.03030303030303,25.54545454545455,34.28787878787879,33,1,1,1,1,1,1),
    @(3,83,"kitchens","living_rooms","distractor","f","stimuli/img_ac0ey.png",86.62222222222222,70.02222222222223,78.32222222222222,45,9,9,9,9,9,9),
    @(4,84,"kitchens","kitchens","target","j","stimuli/img_p659z.png",84.21621621621621,65.37837837837837,74.79729729729729,37,9,9,9,9,9,9),
    @(5,85,"kitchens","kitchens","target","j","stimuli/img_inqod.png",70.84848484848484,50.63636363636363,60.74242424242424,33,5,5,5,5,5,5),
    @(6,86,"kitchens","bedrooms","distractor","f","stimuli/img_dkqas.png",78.57894736842105,57.71052631578947,68.14473684210526,38,7,7,7,7,7,7),
    @(7,87,"kitchens","kitchens","target","j","stimuli/img_aplao.png",64.0909090909091,40.75757575757576,52.42424242424242,33,3,3,3,3,3,3),
    @(8,88,"kitchens","bedrooms","distractor","f","stimuli/img_ybbmx.png",55.24324324324324,36.75675675675676,46,37,3,3,3,3,3,3),
    @(9,89,"kitchens","kitchens","target","j","stimuli/img_yeh72.png",68.66666666666667,45.21212121212121,56.93939393939394,33,4,4,4,4,4,4),
    @(10,90,"kitchens","kitchens","target","j","stimuli/img_ye5sl.png",53.2258064516129,34.45161290322581,43.83870967741936,31,2,2,2,2,2,2),
    @(11,91,"kitchens","kitchens","target","j","stimuli/img_t90e2.png",83.0625,61.96875,72.515625,32,9,9,9,9,9,9),
    @(12,92,"kitchens","kitchens","target","j","stimuli/img_es7o2.png",52.48571428571429,27.54285714285714,40.01428571428572,35,2,2,2,2,2,2),
    @(13,93,"kitchens","bedrooms","distractor","f","stimuli/img_d26ik.png",77.73809523809524,60.66666666666666,69.20238095238095,42,7,7,7,7,7,7),
    @(14,94,"kitchens","bedrooms","distractor","f","stimuli/img_d9ogj.png",76.86842105263158,53.5,65.18421052631578,38,6,6,6,6,6,6),
    @(15,95,"kitchens","living_rooms","distractor","f","stimuli/img_amsgw.png",86.08510638297872,65.95744680851064,76.02127659574468,47,9,9,9,8,9,8),
    @(16,96,"kitchens","kitchens","target","j","stimuli/img_nyv2b.png",11.91176470588235,6.852941176470588,9.382352941176471,34,1,1,1,1,1,1),
    @(17,97,"kitchens","bedrooms","distractor","f","stimuli/img_u1rxv.png",75.63636363636364,54.27272727272727,64.95454545454545,44,6,6,6,6,6,6),
    @(18,98,"kitchens","bedrooms","distractor","f","stimuli/img_n9xll.png",77.14285714285714,59.21428571428572,68.17857142857143,42,7,7,7,7,7,7),
    @(19,99,"kitchens","kitchens","target","j","stimuli/img_eatdk.png",81.40625,61.375,71.390625,32,8,8,8,8,8,8),
    @(20,100,"kitchens","living_rooms","distractor","f","stimuli/img_6zz63.png",87.66666666666667,70.6,79.13333333333333,45,9,10,10,9,9,9),
    @(21,101,"kitchens","kitchens","target","j","stimuli/img_6nbgt.png",78.45161290322581,57.83870967741935,68.14516129032258,31,7,7,7,7,7,7),
    @(22,102,"kitchens","bedrooms","distractor","f","stimuli/img_fea1z.png",79.45945945945945,56.24324324324324,67.85135135135135,37,7,7,7,7,7,7),
    @(23,103,"kitchens","kitchens","target","j","stimuli/img_ce9vx.png",75.9090909090909,57.12121212121212,66.51515151515152,33,7,7,7,7,7,7),
    @(24,104,"kitchens","living_rooms","distractor","f","stimuli/img_pna7l.png",85.53333333333333,67.97777777777777,76.75555555555556,45,9,9,9,9,9,9),
    @(25,105,"kitchens","bedrooms","distractor","f","stimuli/img_iqmdm.png",79.38888888888889,58.36111111111111,68.875,36,7,7,7,7,7,7),
    @(26,106,"kitchens","living_rooms","distractor","f","stimuli/img_yosqb.png",50.88372093023256,30.11627906976744,40.5,43,3,3,3,3,3,3),
    @(27,107,"kitchens","living_rooms","distractor","f","stimuli/img_1zhz6.png",49.02272727272727,32.77272727272727,40.89772727272727,44,3,3,3,3,3,3),
    @(28,108,"kitchens","kitchens","target","j","stimuli/img_d8xbu.png",91.36363636363636,73.18181818181819,82.27272727272728,33,10,10,10,10,10,10),
    @(29,109,"kitchens","kitchens","target","j","stimuli/img_p3hpc.png",72.83333333333333,52.22222222222222,62.52777777777777,36,6,6,6,6,6,6),
    @(30,110,"kitchens","kitchens","target","j","stimuli/img_cnyac.png",69.1470588235294,47.8235294117647,58.48529411764706,34,5,5,5,5,5,5),
    @(31,111,"kitchens","kitchens","target","j","stimuli/img_30vz5.png",86.21212121212122,68.27272727272727,77.24242424242425,33,10,10,10,10,10,10),
    @(32,112,"kitchens","living_rooms","distractor","f","stimuli/img_0jzz7.png",84.85106382978724,68.87234042553192,76.86170212765958,47,9,9,9,9,8,9),
    @(33,113,"kitchens","living_rooms","distractor","f","stimuli/img_wgkqa.png",87.25581395348837,71.13953488372093,79.19767441860465,43,10,10,10,9,9,9),
    @(34,114,"kitchens","kitchens","target","j","stimuli/img_wyl6z.png",59.8235294117647,36.23529411764706,48.02941176470588,34,3,3,3,3,3,3),
    @(35,115,"kitchens","kitchens","target","j","stimuli/img_iyxnj.png",75.30555555555556,54.33333333333334,64.81944444444444,36,6,6,6,6,6,6),
    @(36,116,"kitchens","kitchens","target","j","stimuli/img_njmgp.png",80.48148148148148,58.4074074074074,69.44444444444444,27,8,8,8,8,8,8),
    @(37,117,"kitchens","kitchens","target","j","stimuli/img_cv6mf.png",66.8,42.08,54.44,25,4,4,4,4,4,4),
    @(38,118,"kitchens","bedrooms","distractor","f","stimuli/img_5m6x4.png",80.23076923076923,58.41025641025641,69.32051282051282,39,7,7,7,7,7,7),
    @(39,119,"kitchens","living_rooms","distractor","f","stimuli/img_95hiv.png",84.04545454545455,67.31818181818181,75.68181818181819,44,9,9,9,8,8,8),
    @(40,120,"kitchens","living_rooms","distractor","f","stimuli/img_f63yi.png",85.275,68.475,76.875,40,9,9,9,9,8,9),
    @(41,121,"kitchens","bedrooms","distractor","f","stimuli/img_8dacu.png",76.38461538461539,53.64102564102564,65.01282051282051,39,6,6,6,6,6,6)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value  = $row[1]   # F trial_total
    $ws.Cells.Item($r, 7).Value  = $row[2]   # G target_cat
    $ws.Cells.Item($r, 8).Value  = $row[3]   # H category
    $ws.Cells.Item($r, 9).Value  = $row[4]   # I cond_cat
    $ws.Cells.Item($r, 11).Value = $row[5]   # K correct_answer
    $ws.Cells.Item($r, 12).Value = $row[6]   # L stimulus
    $ws.Cells.Item($r, 13).Value = $row[7]   # M conceptual
    $ws.Cells.Item($r, 14).Value = $row[8]   # N perceptual
    $ws.Cells.Item($r, 15).Value = $row[9]   # O typicality
    $ws.Cells.Item($r, 16).Value = $row[10]  # P n
    $ws.Cells.Item($r, 17).Value = $row[11]  # Q p_typicality
    $ws.Cells.Item($r, 18).Value = $row[12]  # R p_conceptual
    $ws.Cells.Item($r, 19).Value = $row[13]  # S p_perceptual
    $ws.Cells.Item($r, 20).Value = $row[14]  # T r_typicality
    $ws.Cells.Item($r, 21).Value = $row[15]  # U r_conceptual
    $ws.Cells.Item($r, 22).Value = $row[16]  # V r_perceptual
}
